$d = $word.ActiveDocument

# Locate the paragraph containing the "LOQ4083" requirement line, then
# remove the three paragraphs that follow it:
#   1) the blank paragraph right after it
#   2) "Ver no Jupiter Salvar em pdf Salvar em docx"
#   3) the "© 2020 . Contact: ..." paragraph
# The blank paragraph that precedes the page-break paragraph at the very
# end of the document must be left untouched.

$anchorText = "LOQ4083: Fen"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$anchorText*") {
        $target = $i
    }
}

if ($target -ne $null) {
    $firstToDelete = $target + 1
    $lastToDelete = $target + 3

    $startRange = $d.Paragraphs.Item($firstToDelete).Range
    $endRange = $d.Paragraphs.Item($lastToDelete).Range

    $deleteRange = $d.Range($startRange.Start, $endRange.End)
    $deleteRange.Delete()
}
